{"js": "const replacements = [\n  [\"17+15=32\", \"49+34=83\"],\n  [\"63+29=92\", \"27+47=74\"],\n  [\"5+28=33\", \"39+39=78\"],\n  [\"51-22=29\", \"85-47=38\"],\n  [\"86+9=95\", \"40-23=17\"],\n  [\"18+13=31\", \"16+8=24\"],\n  [\"77-9=68\", \"17+49=66\"],\n  [\"84-55=29\", \"51-7=44\"],\n  [\"27+68=95\", \"75+18=93\"],\n  [\"15+59=74\", \"48+13=61\"],\n  [\"3+18=21\", \"28+17=45\"],\n  [\"94-9=85\", \"17+77=94\"],\n  [\"73-54=19\", \"60-1=59\"],\n  [\"80-41=39\", \"58+23=81\"],\n  [\"53-46=7\", \"96-48=48\"],\n  [\"72-59=13\", \"53-14=39\"],\n  [\"70-39=31\", \"71-47=24\"],\n  [\"37-8=29\", \"57+39=96\"],\n  [\"55-36=19\", \"8+29=37\"],\n  [\"7+36=43\", \"43-34=9\"],\n  [\"86-29=57\", \"80-51=29\"],\n  [\"39+17=56\", \"70-63=7\"],\n  [\"82-16=66\", \"46+46=92\"],\n  [\"7+38=45\", \"28+34=62\"],\n  [\"91-23=68\", \"13-7=6\"],\n  [\"18+16=34\", \"83-15=68\"],\n  [\"73-56=17\", \"90-73=17\"],\n  [\"32-29=3\", \"91-6=85\"],\n  [\"85-77=8\", \"32-25=7\"],\n  [\"18+7=25\", \"29+9=38\"],\n  [\"6+58=64\", \"56+15=71\"],\n  [\"88+6=94\", \"80-25=55\"],\n  [\"70-29=41\", \"67-29=38\"],\n  [\"56+7=63\", \"73+18=91\"],\n  [\"46+36=82\", \"20-5=15\"],\n  [\"52-17=35\", \"9+66=75\"],\n  [\"9+55=64\", \"8+75=83\"],\n  [\"62-54=8\", \"64-56=8\"],\n  [\"85-56=29\", \"76-37=39\"],\n  [\"63-15=48\", \"16+16=32\"],\n  [\"68+5=73\", \"9+4=13\"],\n  [\"56+38=94\", \"91-79=12\"],\n  [\"8+19=27\", \"17+7=24\"],\n  [\"62-47=15\", \"13-7=6\"],\n  [\"41-17=24\", \"83-68=15\"],\n  [\"39+54=93\", \"82-25=57\"],\n  [\"65-48=17\", \"92-83=9\"],\n  [\"35-16=19\", \"81-19=62\"],\n  [\"62-28=34\", \"33-16=17\"],\n  [\"85-58=27\", \"62-19=43\"],\n  [\"86+7=93\", \"98-49=49\"],\n  [\"29+12=41\", \"73-27=46\"],\n  [\"47+5=52\", \"4+89=93\"],\n  [\"91-36=55\", \"60-48=12\"],\n  [\"38+45=83\", \"92-45=47\"],\n  [\"80-9=71\", \"43-38=5\"],\n  [\"22+39=61\", \"54-47=7\"],\n  [\"36-27=9\", \"52-9=43\"],\n  [\"50-35=15\", \"98-79=19\"],\n  [\"76-39=37\", \"40-16=24\"],\n  [\"77-68=9\", \"35-26=9\"],\n  [\"63-16=47\", \"59+29=88\"],\n  [\"18+24=42\", \"69+5=74\"],\n  [\"45+37=82\", \"68+17=85\"],\n  [\"36+35=71\", \"65-17=48\"],\n  [\"19+63=82\", \"38+43=81\"],\n  [\"8+28=36\", \"15+47=62\"],\n  [\"51-35=16\", \"15+77=92\"],\n  [\"6+36=42\", \"76+5=81\"],\n  [\"18+68=86\", \"70-65=5\"],\n  [\"61-47=14\", \"95-58=37\"],\n  [\"30-22=8\", \"50-48=2\"],\n  [\"17+45=62\", \"40-2=38\"],\n  [\"95-39=56\", \"80-7=73\"],\n  [\"75-18=57\", \"31-17=14\"],\n  [\"18+29=47\", \"47+45=92\"],\n  [\"80-18=62\", \"42-5=37\"],\n  [\"70-1=69\", \"69+5=74\"],\n  [\"68+13=81\", \"96-29=67\"],\n  [\"28+14=42\", \"25+16=41\"],\n  [\"34-27=7\", \"76-49=27\"],\n  [\"47+16=63\", \"5+38=43\"],\n  [\"6+29=35\", \"91-58=33\"],\n  [\"63+18=81\", \"57+28=85\"],\n  [\"8+17=25\", \"83-5=78\"],\n  [\"84-78=6\", \"58+39=97\"],\n  [\"86+5=91\", \"75-6=69\"],\n  [\"39+54=93\", \"70-6=64\"],\n  [\"14+47=61\", \"18+55=73\"],\n  [\"10-1=9\", \"71-66=5\"],\n  [\"55-9=46\", \"30-9=21\"],\n  [\"17+27=44\", \"34-26=8\"],\n  [\"32+49=81\", \"84-65=19\"],\n  [\"49+16=65\", \"57+5=62\"],\n  [\"20-1=19\", \"23+49=72\"],\n  [\"83-5=78\", \"90-23=67\"],\n  [\"33-25=8\", \"27+37=64\"],\n  [\"19+43=62\", \"69+14=83\"],\n  [\"8+33=41\", \"65-47=18\"],\n  [\"75-19=56\", \"42-8=34\"]\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Skip the first paragraph (title/date); the remaining paragraphs are the\n// 100 math-expression cells, in document order, one per replacement.\nconst mathParagraphs = paragraphs.items.slice(1);\n\nif (mathParagraphs.length !== replacements.length) {\n  throw new Error(\n    \"Expected \" + replacements.length + \" math paragraphs, found \" + mathParagraphs.length\n  );\n}\n\nfor (let i = 0; i < replacements.length; i++) {\n  const newText = replacements[i][1];\n  mathParagraphs[i].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document body is a single 20x5 table of arithmetic problems\n# (\"a+b=c\" / \"a-b=c\"), one per cell, plus a title paragraph above it.\n# The replacements below are the old/new cell text values, in table\n# (row-major) order: row 1 col 1..5, row 2 col 1..5, etc.\n$replacements = @(\n    @('17+15=32', '49+34=83'),\n    @('63+29=92', '27+47=74'),\n    @('5+28=33', '39+39=78'),\n    @('51-22=29', '85-47=38'),\n    @('86+9=95', '40-23=17'),\n    @('18+13=31', '16+8=24'),\n    @('77-9=68', '17+49=66'),\n    @('84-55=29', '51-7=44'),\n    @('27+68=95', '75+18=93'),\n    @('15+59=74', '48+13=61'),\n    @('3+18=21', '28+17=45'),\n    @('94-9=85', '17+77=94'),\n    @('73-54=19', '60-1=59'),\n    @('80-41=39', '58+23=81'),\n    @('53-46=7', '96-48=48'),\n    @('72-59=13', '53-14=39'),\n    @('70-39=31', '71-47=24'),\n    @('37-8=29', '57+39=96'),\n    @('55-36=19', '8+29=37'),\n    @('7+36=43', '43-34=9'),\n    @('86-29=57', '80-51=29'),\n    @('39+17=56', '70-63=7'),\n    @('82-16=66', '46+46=92'),\n    @('7+38=45', '28+34=62'),\n    @('91-23=68', '13-7=6'),\n    @('18+16=34', '83-15=68'),\n    @('73-56=17', '90-73=17'),\n    @('32-29=3', '91-6=85'),\n    @('85-77=8', '32-25=7'),\n    @('18+7=25', '29+9=38'),\n    @('6+58=64', '56+15=71'),\n    @('88+6=94', '80-25=55'),\n    @('70-29=41', '67-29=38'),\n    @('56+7=63', '73+18=91'),\n    @('46+36=82', '20-5=15'),\n    @('52-17=35', '9+66=75'),\n    @('9+55=64', '8+75=83'),\n    @('62-54=8', '64-56=8'),\n    @('85-56=29', '76-37=39'),\n    @('63-15=48', '16+16=32'),\n    @('68+5=73', '9+4=13'),\n    @('56+38=94', '91-79=12'),\n    @('8+19=27', '17+7=24'),\n    @('62-47=15', '13-7=6'),\n    @('41-17=24', '83-68=15'),\n    @('39+54=93', '82-25=57'),\n    @('65-48=17', '92-83=9'),\n    @('35-16=19', '81-19=62'),\n    @('62-28=34', '33-16=17'),\n    @('85-58=27', '62-19=43'),\n    @('86+7=93', '98-49=49'),\n    @('29+12=41', '73-27=46'),\n    @('47+5=52', '4+89=93'),\n    @('91-36=55', '60-48=12'),\n    @('38+45=83', '92-45=47'),\n    @('80-9=71', '43-38=5'),\n    @('22+39=61', '54-47=7'),\n    @('36-27=9', '52-9=43'),\n    @('50-35=15', '98-79=19'),\n    @('76-39=37', '40-16=24'),\n    @('77-68=9', '35-26=9'),\n    @('63-16=47', '59+29=88'),\n    @('18+24=42', '69+5=74'),\n    @('45+37=82', '68+17=85'),\n    @('36+35=71', '65-17=48'),\n    @('19+63=82', '38+43=81'),\n    @('8+28=36', '15+47=62'),\n    @('51-35=16', '15+77=92'),\n    @('6+36=42', '76+5=81'),\n    @('18+68=86', '70-65=5'),\n    @('61-47=14', '95-58=37'),\n    @('30-22=8', '50-48=2'),\n    @('17+45=62', '40-2=38'),\n    @('95-39=56', '80-7=73'),\n    @('75-18=57', '31-17=14'),\n    @('18+29=47', '47+45=92'),\n    @('80-18=62', '42-5=37'),\n    @('70-1=69', '69+5=74'),\n    @('68+13=81', '96-29=67'),\n    @('28+14=42', '25+16=41'),\n    @('34-27=7', '76-49=27'),\n    @('47+16=63', '5+38=43'),\n    @('6+29=35', '91-58=33'),\n    @('63+18=81', '57+28=85'),\n    @('8+17=25', '83-5=78'),\n    @('84-78=6', '58+39=97'),\n    @('86+5=91', '75-6=69'),\n    @('39+54=93', '70-6=64'),\n    @('14+47=61', '18+55=73'),\n    @('10-1=9', '71-66=5'),\n    @('55-9=46', '30-9=21'),\n    @('17+27=44', '34-26=8'),\n    @('32+49=81', '84-65=19'),\n    @('49+16=65', '57+5=62'),\n    @('20-1=19', '23+49=72'),\n    @('83-5=78', '90-23=67'),\n    @('33-25=8', '27+37=64'),\n    @('19+43=62', '69+14=83'),\n    @('8+33=41', '65-47=18'),\n    @('75-19=56', '42-8=34')\n)\n\n$table = $d.Tables(1)\n$rows = $table.Rows.Count\n$cols = $table.Columns.Count\n\nif (($rows * $cols) -ne $replacements.Count) {\n    throw \"Expected $($replacements.Count) cells, found $($rows * $cols)\"\n}\n\n$i = 0\nfor ($row = 1; $row -le $rows; $row++) {\n    for ($col = 1; $col -le $cols; $col++) {\n        $pair = $replacements[$i]\n        $oldText = $pair[0]\n        $newText = $pair[1]\n\n        $cell = $table.Cell($row, $col)\n        $cellRange = $cell.Range\n        # Exclude the trailing paragraph mark / end-of-cell marker so only\n        # the visible text content is replaced.\n        $textRange = $d.Range($cellRange.Start, $cellRange.End - 1)\n        $textRange.Text = $newText\n\n        $i++\n    }\n}\n"}
